$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9

# Plain numeric columns
$ws.Cells.Item($row, 1).Value = 131146223    # A - Id
$ws.Cells.Item($row, 2).Value = 57881        # B - Taxonsorteringsordning
$ws.Cells.Item($row, 5).Value = 100049       # E - TaxonId
$ws.Cells.Item($row, 17).Value = 592302      # Q - Ost
$ws.Cells.Item($row, 18).Value = 6320843     # R - Nord
$ws.Cells.Item($row, 19).Value = 25          # S - Noggrannhet

# Plain text columns (not numeric-looking, stay text naturally)
$ws.Cells.Item($row, 4).Value = "NT"                                   # D - Rodlistade
$ws.Cells.Item($row, 6).Value = "Spillkråka"                           # F - Artnamn
$ws.Cells.Item($row, 7).Value = "Dryocopus martius"                    # G - Vetenskapligt namn
$ws.Cells.Item($row, 8).Value = "(Linnaeus, 1758)"                     # H - Auktor
$ws.Cells.Item($row, 13).Value = "lockläte, övriga läten"              # M - Aktivitet
$ws.Cells.Item($row, 16).Value = "Korshamn, Björnö, Mönsterås, Sm"     # P - Lokalnamn
$ws.Cells.Item($row, 20).Value = "Kalmar"                              # T - Lan
$ws.Cells.Item($row, 21).Value = "Mönsterås"                           # U - Kommun
$ws.Cells.Item($row, 22).Value = "Småland"                             # V - Provins
$ws.Cells.Item($row, 23).Value = "Mönsterås"                           # W - Socken
$ws.Cells.Item($row, 26).Value = "13:00"                               # Z - Starttid
$ws.Cells.Item($row, 28).Value = "15:30"                               # AB - Sluttid
$ws.Cells.Item($row, 49).Value = "Jan Brenander"                       # AW - Rapportor
$ws.Cells.Item($row, 50).Value = "Jan Brenander"                       # AX - Observatorer

# Text columns that look numeric / date-like: force text storage with a
# leading apostrophe so Excel keeps them as strings instead of converting
# them to numbers or date serials.
$ws.Cells.Item($row, 9).Value = "'1"              # I - Antal
$ws.Cells.Item($row, 25).Value = "'2026-02-13"    # Y - Startdatum
$ws.Cells.Item($row, 27).Value = "'2026-02-13"    # AA - Slutdatum

# Boolean columns
$ws.Cells.Item($row, 30).Value = $false   # AD - Ej återfunnen
$ws.Cells.Item($row, 31).Value = $false   # AE - Osaker artbestamning
$ws.Cells.Item($row, 33).Value = $false   # AG - Ospontan
